$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data Insights")

# Insert a new row at row 10, pushing existing rows 10-20 down to 11-21
$ws.Rows.Item(10).Insert()

# Add the new drill-through note in column B of the newly inserted row
$ws.Range("B10").Value = "In 2015 Quarter 4, there is a discrepancy, December_Sales > October_Sales but December_Profit < October_Profit. Need to investigate."

# Re-select the cell Excel ends up with after this edit
$ws.Range("B14").Select()

# Column B needs to widen to fit the new, longer text (bestFit / AutoFit)
$ws.Columns.Item(2).ColumnWidth = 120.6
